$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value = 15
$ws1.Range("F3").Value = 25
$ws1.Range("F5").Value = 5034
$ws1.Range("F6").Value = 5034
$ws1.Range("F7").Value = 78
$ws1.Range("F9").Value = 497
$ws1.Range("F11").Value = 1139
$ws1.Range("F12").Value = 691
$ws1.Range("F13").Value = 4848
$ws1.Range("F14").Value = 20
$ws1.Range("F15").Value = 49
$ws1.Range("F16").Value = 67
$ws1.Range("F17").Value = 201
$ws1.Range("F18").Value = 207
$ws1.Range("F19").Value = 93
$ws1.Range("F20").Value = 240
$ws1.Range("F21").Value = 3720
$ws1.Range("F23").Value = 35
$ws1.Range("F24").Value = 3568
$ws1.Range("F25").Value = 165
$ws1.Range("F26").Value = 157
$ws1.Range("F28").Value = 190
$ws1.Range("F29").Value = 227
$ws1.Range("F30").Value = 197
$ws1.Range("F35").Value = 134
$ws1.Range("F36").Value = 6272
$ws1.Range("F37").Value = 991
$ws1.Range("F38").Value = 473
$ws1.Range("F39").Value = 93
$ws1.Range("F42").Value = 1290
$ws1.Range("F44").Value = 619
$ws1.Range("F46").Value = 2166
$ws1.Range("F49").Value = 751
$ws1.Range("F50").Value = 894

# Sheet 2: 演出
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F3").Value = 18
$ws2.Range("F20").Value = 46
$ws2.Range("F23").Value = 793

# Sheet 4: 全部类型
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F4").Value = 18
$ws4.Range("F5").Value = 15
$ws4.Range("F8").Value = 25
$ws4.Range("F10").Value = 5034
$ws4.Range("F11").Value = 5034
$ws4.Range("F12").Value = 78
$ws4.Range("F16").Value = 691
$ws4.Range("F17").Value = 4848
$ws4.Range("F18").Value = 20
$ws4.Range("F19").Value = 49
$ws4.Range("F20").Value = 67
$ws4.Range("F21").Value = 201
$ws4.Range("F22").Value = 93
$ws4.Range("F23").Value = 240
$ws4.Range("F24").Value = 3568
$ws4.Range("F25").Value = 165
$ws4.Range("F26").Value = 157
$ws4.Range("F27").Value = 190
$ws4.Range("F28").Value = 227
$ws4.Range("F29").Value = 197
$ws4.Range("F34").Value = 134
$ws4.Range("F36").Value = 6272
$ws4.Range("F37").Value = 991
$ws4.Range("F38").Value = 93
$ws4.Range("F40").Value = 1290
$ws4.Range("F42").Value = 619
$ws4.Range("F44").Value = 2166
$ws4.Range("F48").Value = 751
$ws4.Range("F49").Value = 894
